$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 25.55000000000055
$ws.Range("G2").Value = [double]"1.362692110262742e-07"
$ws.Range("H2").Value = [double]"2.637387653543786e-06"
$ws.Range("K2").Value = 5.550980196642649
$ws.Range("L2").Value = "[3.2386660630068356, 7.863294330278462]"
$ws.Range("M2").Value = [double]"3.316484054716184e-06"
$ws.Range("N2").Value = [double]"3.316484054716184e-06"
$ws.Range("O2").Value = -1.383684452031541
$ws.Range("P2").Value = "[-1.8491055858966954, -0.9182633181663862]"
$ws.Range("Q2").Value = [double]"1.086076317768914e-08"
$ws.Range("R2").Value = [double]"1.086076317768914e-08"
$ws.Range("S2").Value = 10.70317608807971
$ws.Range("T2").Value = "[9.360875628743212, 12.04547654741621]"
$ws.Range("W2").Value = 5.626626626626749
$ws.Range("X2").Value = 3.734034034034114
$ws.Range("Y2").Value = 7.519219219219383

# Row 3 updates
$ws.Range("E3").Value = 24.60000000000041
$ws.Range("G3").Value = [double]"1.692220250482102e-06"
$ws.Range("H3").Value = [double]"8.891491867199399e-06"
$ws.Range("I3").Value = [double]"4.485301019485632e-14"
$ws.Range("K3").Value = 5.507892094037089
$ws.Range("L3").Value = "[3.2442875762070553, 7.771496611867123]"
$ws.Range("M3").Value = [double]"2.693249301799483e-06"
$ws.Range("N3").Value = [double]"3.316484054716184e-06"
$ws.Range("O3").Value = -3.107000542289004
$ws.Range("P3").Value = "[-3.6478953735376978, -2.5661057110403105]"
$ws.Range("S3").Value = 9.679426197961195
$ws.Range("T3").Value = "[8.212458606273167, 11.146393789649224]"
$ws.Range("W3").Value = 12.16456456456477
$ws.Range("X3").Value = 10.04684684684701
$ws.Range("Y3").Value = 14.28228228228252
